# Auto-generated: append new daily snapshot rows to each portfolio sheet
# (大智投资组合收益, 大成投资组合收益, 我的投资组合收益) matching the
# upstream commit's newly-appended rows. Numeric-looking text (stock
# codes / timestamps) is written via a quote-prefixed Formula so Excel
# keeps it as text (preserving leading zeros) instead of coercing it to
# a number; the temporary quotePrefix style is then reset to Normal so
# no stray formatting is left behind.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(72, 1).Formula = '''大智 (稳健智远)'
$ws.Cells.Item(72, 1).Style = "Normal"
$ws.Cells.Item(72, 2).Formula = '''000333'
$ws.Cells.Item(72, 2).Style = "Normal"
$ws.Cells.Item(72, 3).Formula = '''美的集团'
$ws.Cells.Item(72, 3).Style = "Normal"
$ws.Cells.Item(72, 4).Value = 3.03
$ws.Cells.Item(72, 5).Value = 42.89719527444501
$ws.Cells.Item(72, 6).Value = 72.06
$ws.Cells.Item(72, 7).Value = 3091.171891476507
$ws.Cells.Item(72, 8).Value = 102088.3450563826
$ws.Cells.Item(72, 9).Formula = '''202506171600'
$ws.Cells.Item(72, 9).Style = "Normal"
$ws.Cells.Item(73, 1).Formula = '''大智 (稳健智远)'
$ws.Cells.Item(73, 1).Style = "Normal"
$ws.Cells.Item(73, 2).Formula = '''510050'
$ws.Cells.Item(73, 2).Style = "Normal"
$ws.Cells.Item(73, 3).Formula = '''上证50ETF'
$ws.Cells.Item(73, 3).Style = "Normal"
$ws.Cells.Item(73, 4).Value = 4.99
$ws.Cells.Item(73, 5).Value = 1852.638869852698
$ws.Cells.Item(73, 6).Value = 2.75
$ws.Cells.Item(73, 7).Value = 5094.75689209492
$ws.Cells.Item(73, 8).Value = 102088.3450563826
$ws.Cells.Item(73, 9).Formula = '''202506171600'
$ws.Cells.Item(73, 9).Style = "Normal"
$ws.Cells.Item(74, 1).Formula = '''大智 (稳健智远)'
$ws.Cells.Item(74, 1).Style = "Normal"
$ws.Cells.Item(74, 2).Formula = '''510300'
$ws.Cells.Item(74, 2).Style = "Normal"
$ws.Cells.Item(74, 3).Formula = '''沪深300ETF'
$ws.Cells.Item(74, 3).Style = "Normal"
$ws.Cells.Item(74, 4).Value = 4.99
$ws.Cells.Item(74, 5).Value = 1276.881426590205
$ws.Cells.Item(74, 6).Value = 3.99
$ws.Cells.Item(74, 7).Value = 5094.756892094918
$ws.Cells.Item(74, 8).Value = 102088.3450563826
$ws.Cells.Item(74, 9).Formula = '''202506171600'
$ws.Cells.Item(74, 9).Style = "Normal"
$ws.Cells.Item(75, 1).Formula = '''大智 (稳健智远)'
$ws.Cells.Item(75, 1).Style = "Normal"
$ws.Cells.Item(75, 2).Formula = '''518880'
$ws.Cells.Item(75, 2).Style = "Normal"
$ws.Cells.Item(75, 3).Formula = '''黄金ETF'
$ws.Cells.Item(75, 3).Style = "Normal"
$ws.Cells.Item(75, 4).Value = 4.95
$ws.Cells.Item(75, 5).Value = 673.9096418114974
$ws.Cells.Item(75, 6).Value = 7.5
$ws.Cells.Item(75, 7).Value = 5054.322313586231
$ws.Cells.Item(75, 8).Value = 102088.3450563826
$ws.Cells.Item(75, 9).Formula = '''202506171600'
$ws.Cells.Item(75, 9).Style = "Normal"
$ws.Cells.Item(76, 1).Formula = '''大智 (稳健智远)'
$ws.Cells.Item(76, 1).Style = "Normal"
$ws.Cells.Item(76, 2).Formula = '''600085'
$ws.Cells.Item(76, 2).Style = "Normal"
$ws.Cells.Item(76, 3).Formula = '''同仁堂'
$ws.Cells.Item(76, 3).Style = "Normal"
$ws.Cells.Item(76, 4).Value = 1.99
$ws.Cells.Item(76, 5).Value = 56.06334956913254
$ws.Cells.Item(76, 6).Value = 36.19
$ws.Cells.Item(76, 7).Value = 2028.932620906907
$ws.Cells.Item(76, 8).Value = 102088.3450563826
$ws.Cells.Item(76, 9).Formula = '''202506171600'
$ws.Cells.Item(76, 9).Style = "Normal"
$ws.Cells.Item(77, 1).Formula = '''大智 (稳健智远)'
$ws.Cells.Item(77, 1).Style = "Normal"
$ws.Cells.Item(77, 2).Formula = '''600900'
$ws.Cells.Item(77, 2).Style = "Normal"
$ws.Cells.Item(77, 3).Formula = '''长江电力'
$ws.Cells.Item(77, 3).Style = "Normal"
$ws.Cells.Item(77, 4).Value = 20.1
$ws.Cells.Item(77, 5).Value = 669.0422707938175
$ws.Cells.Item(77, 6).Value = 30.67
$ws.Cells.Item(77, 7).Value = 20519.52644524638
$ws.Cells.Item(77, 8).Value = 102088.3450563826
$ws.Cells.Item(77, 9).Formula = '''202506171600'
$ws.Cells.Item(77, 9).Style = "Normal"
$ws.Cells.Item(78, 1).Formula = '''大智 (稳健智远)'
$ws.Cells.Item(78, 1).Style = "Normal"
$ws.Cells.Item(78, 2).Formula = '''600989'
$ws.Cells.Item(78, 2).Style = "Normal"
$ws.Cells.Item(78, 3).Formula = '''宝丰能源'
$ws.Cells.Item(78, 3).Style = "Normal"
$ws.Cells.Item(78, 4).Value = 4.97
$ws.Cells.Item(78, 5).Value = 308.5861230826723
$ws.Cells.Item(78, 6).Value = 16.43
$ws.Cells.Item(78, 7).Value = 5070.070002248306
$ws.Cells.Item(78, 8).Value = 102088.3450563826
$ws.Cells.Item(78, 9).Formula = '''202506171600'
$ws.Cells.Item(78, 9).Style = "Normal"
$ws.Cells.Item(79, 1).Formula = '''大智 (稳健智远)'
$ws.Cells.Item(79, 1).Style = "Normal"
$ws.Cells.Item(79, 2).Formula = '''601899'
$ws.Cells.Item(79, 2).Style = "Normal"
$ws.Cells.Item(79, 3).Formula = '''XD紫金矿'
$ws.Cells.Item(79, 3).Style = "Normal"
$ws.Cells.Item(79, 4).Value = 10.01
$ws.Cells.Item(79, 5).Value = 541.4194359293219
$ws.Cells.Item(79, 6).Value = 18.87
$ws.Cells.Item(79, 7).Value = 10216.58475598631
$ws.Cells.Item(79, 8).Value = 102088.3450563826
$ws.Cells.Item(79, 9).Formula = '''202506171600'
$ws.Cells.Item(79, 9).Style = "Normal"
$ws.Cells.Item(80, 1).Formula = '''大智 (稳健智远)'
$ws.Cells.Item(80, 1).Style = "Normal"
$ws.Cells.Item(80, 2).Formula = '''HK02899'
$ws.Cells.Item(80, 2).Style = "Normal"
$ws.Cells.Item(80, 3).Formula = '''紫金矿业'
$ws.Cells.Item(80, 3).Style = "Normal"
$ws.Cells.Item(80, 4).Value = 10.06
$ws.Cells.Item(80, 5).Value = 521.4694874201556
$ws.Cells.Item(80, 6).Value = 19.7
$ws.Cells.Item(80, 7).Value = 10272.94890217706
$ws.Cells.Item(80, 8).Value = 102088.3450563826
$ws.Cells.Item(80, 9).Formula = '''202506171600'
$ws.Cells.Item(80, 9).Style = "Normal"
$ws.Cells.Item(81, 1).Formula = '''大智 (稳健智远)'
$ws.Cells.Item(81, 1).Style = "Normal"
$ws.Cells.Item(81, 2).Formula = '''HK06881'
$ws.Cells.Item(81, 2).Style = "Normal"
$ws.Cells.Item(81, 3).Formula = '''中国银河'
$ws.Cells.Item(81, 3).Style = "Normal"
$ws.Cells.Item(81, 4).Value = 4.97
$ws.Cells.Item(81, 5).Value = 600.7968033130801
$ws.Cells.Item(81, 6).Value = 8.45
$ws.Cells.Item(81, 7).Value = 5076.732987995526
$ws.Cells.Item(81, 8).Value = 102088.3450563826
$ws.Cells.Item(81, 9).Formula = '''202506171600'
$ws.Cells.Item(81, 9).Style = "Normal"
$ws.Cells.Item(82, 1).Formula = '''大智 (稳健智远)'
$ws.Cells.Item(82, 1).Style = "Normal"
$ws.Cells.Item(82, 2).Formula = '''100000'
$ws.Cells.Item(82, 2).Style = "Normal"
$ws.Cells.Item(82, 3).Formula = '''现金'
$ws.Cells.Item(82, 3).Style = "Normal"
$ws.Cells.Item(82, 4).Value = 29.94
$ws.Cells.Item(82, 5).Value = 30568.54135256952
$ws.Cells.Item(82, 6).Value = 1
$ws.Cells.Item(82, 7).Value = 30568.54135256952
$ws.Cells.Item(82, 8).Value = 102088.3450563826
$ws.Cells.Item(82, 9).Formula = '''202506171600'
$ws.Cells.Item(82, 9).Style = "Normal"

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(51, 1).Formula = '''大成 (锐进先锋)'
$ws.Cells.Item(51, 1).Style = "Normal"
$ws.Cells.Item(51, 2).Formula = '''000725'
$ws.Cells.Item(51, 2).Style = "Normal"
$ws.Cells.Item(51, 3).Formula = '''京东方A'
$ws.Cells.Item(51, 3).Style = "Normal"
$ws.Cells.Item(51, 4).Value = 5.12
$ws.Cells.Item(51, 5).Value = 1264.047242468624
$ws.Cells.Item(51, 6).Value = 3.92
$ws.Cells.Item(51, 7).Value = 4955.065190477007
$ws.Cells.Item(51, 8).Value = 96797.0862179527
$ws.Cells.Item(51, 9).Formula = '''202506171600'
$ws.Cells.Item(51, 9).Style = "Normal"
$ws.Cells.Item(52, 1).Formula = '''大成 (锐进先锋)'
$ws.Cells.Item(52, 1).Style = "Normal"
$ws.Cells.Item(52, 2).Formula = '''159781'
$ws.Cells.Item(52, 2).Style = "Normal"
$ws.Cells.Item(52, 3).Formula = '''科创创业ETF'
$ws.Cells.Item(52, 3).Style = "Normal"
$ws.Cells.Item(52, 4).Value = 5.08
$ws.Cells.Item(52, 5).Value = 9277.629760760277
$ws.Cells.Item(52, 6).Value = 0.53
$ws.Cells.Item(52, 7).Value = 4917.143773202947
$ws.Cells.Item(52, 8).Value = 96797.0862179527
$ws.Cells.Item(52, 9).Formula = '''202506171600'
$ws.Cells.Item(52, 9).Style = "Normal"
$ws.Cells.Item(53, 1).Formula = '''大成 (锐进先锋)'
$ws.Cells.Item(53, 1).Style = "Normal"
$ws.Cells.Item(53, 2).Formula = '''513100'
$ws.Cells.Item(53, 2).Style = "Normal"
$ws.Cells.Item(53, 3).Formula = '''纳指ETF'
$ws.Cells.Item(53, 3).Style = "Normal"
$ws.Cells.Item(53, 4).Value = 5.08
$ws.Cells.Item(53, 5).Value = 3131.938709046463
$ws.Cells.Item(53, 6).Value = 1.57
$ws.Cells.Item(53, 7).Value = 4917.143773202947
$ws.Cells.Item(53, 8).Value = 96797.0862179527
$ws.Cells.Item(53, 9).Formula = '''202506171600'
$ws.Cells.Item(53, 9).Style = "Normal"
$ws.Cells.Item(54, 1).Formula = '''大成 (锐进先锋)'
$ws.Cells.Item(54, 1).Style = "Normal"
$ws.Cells.Item(54, 2).Formula = '''513290'
$ws.Cells.Item(54, 2).Style = "Normal"
$ws.Cells.Item(54, 3).Formula = '''纳指生物科技ETF'
$ws.Cells.Item(54, 3).Style = "Normal"
$ws.Cells.Item(54, 4).Value = 1.01
$ws.Cells.Item(54, 5).Value = 870.2909333102562
$ws.Cells.Item(54, 6).Value = 1.12
$ws.Cells.Item(54, 7).Value = 974.725845307487
$ws.Cells.Item(54, 8).Value = 96797.0862179527
$ws.Cells.Item(54, 9).Formula = '''202506171600'
$ws.Cells.Item(54, 9).Style = "Normal"
$ws.Cells.Item(55, 1).Formula = '''大成 (锐进先锋)'
$ws.Cells.Item(55, 1).Style = "Normal"
$ws.Cells.Item(55, 2).Formula = '''603119'
$ws.Cells.Item(55, 2).Style = "Normal"
$ws.Cells.Item(55, 3).Formula = '''浙江荣泰'
$ws.Cells.Item(55, 3).Style = "Normal"
$ws.Cells.Item(55, 4).Value = 44.72
$ws.Cells.Item(55, 5).Value = 1069.978093782073
$ws.Cells.Item(55, 6).Value = 40.46
$ws.Cells.Item(55, 7).Value = 43291.31367442268
$ws.Cells.Item(55, 8).Value = 96797.0862179527
$ws.Cells.Item(55, 9).Formula = '''202506171600'
$ws.Cells.Item(55, 9).Style = "Normal"
$ws.Cells.Item(56, 1).Formula = '''大成 (锐进先锋)'
$ws.Cells.Item(56, 1).Style = "Normal"
$ws.Cells.Item(56, 2).Formula = '''688290'
$ws.Cells.Item(56, 2).Style = "Normal"
$ws.Cells.Item(56, 3).Formula = '''景业智能'
$ws.Cells.Item(56, 3).Style = "Normal"
$ws.Cells.Item(56, 4).Value = 8.51
$ws.Cells.Item(56, 5).Value = 147.121987895035
$ws.Cells.Item(56, 6).Value = 56
$ws.Cells.Item(56, 7).Value = 8238.83132212196
$ws.Cells.Item(56, 8).Value = 96797.0862179527
$ws.Cells.Item(56, 9).Formula = '''202506171600'
$ws.Cells.Item(56, 9).Style = "Normal"
$ws.Cells.Item(57, 1).Formula = '''大成 (锐进先锋)'
$ws.Cells.Item(57, 1).Style = "Normal"
$ws.Cells.Item(57, 2).Formula = '''100000'
$ws.Cells.Item(57, 2).Style = "Normal"
$ws.Cells.Item(57, 3).Formula = '''现金'
$ws.Cells.Item(57, 3).Style = "Normal"
$ws.Cells.Item(57, 4).Value = 30.48
$ws.Cells.Item(57, 5).Value = 29502.86263921768
$ws.Cells.Item(57, 6).Value = 1
$ws.Cells.Item(57, 7).Value = 29502.86263921768
$ws.Cells.Item(57, 8).Value = 96797.0862179527
$ws.Cells.Item(57, 9).Formula = '''202506171600'
$ws.Cells.Item(57, 9).Style = "Normal"

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(106, 1).Formula = '''范式进化投资组合'
$ws.Cells.Item(106, 1).Style = "Normal"
$ws.Cells.Item(106, 2).Formula = '''000333'
$ws.Cells.Item(106, 2).Style = "Normal"
$ws.Cells.Item(106, 3).Formula = '''美的集团'
$ws.Cells.Item(106, 3).Style = "Normal"
$ws.Cells.Item(106, 4).Value = 1.01
$ws.Cells.Item(106, 5).Value = 14.02515882310653
$ws.Cells.Item(106, 6).Value = 72.06
$ws.Cells.Item(106, 7).Value = 1010.652944793057
$ws.Cells.Item(106, 8).Value = 99955.80999104423
$ws.Cells.Item(106, 9).Formula = '''202506171600'
$ws.Cells.Item(106, 9).Style = "Normal"
$ws.Cells.Item(107, 1).Formula = '''范式进化投资组合'
$ws.Cells.Item(107, 1).Style = "Normal"
$ws.Cells.Item(107, 2).Formula = '''000725'
$ws.Cells.Item(107, 2).Style = "Normal"
$ws.Cells.Item(107, 3).Formula = '''京东方A'
$ws.Cells.Item(107, 3).Style = "Normal"
$ws.Cells.Item(107, 4).Value = 5.04
$ws.Cells.Item(107, 5).Value = 1284.618017653691
$ws.Cells.Item(107, 6).Value = 3.92
$ws.Cells.Item(107, 7).Value = 5035.702629202468
$ws.Cells.Item(107, 8).Value = 99955.80999104423
$ws.Cells.Item(107, 9).Formula = '''202506171600'
$ws.Cells.Item(107, 9).Style = "Normal"
$ws.Cells.Item(108, 1).Formula = '''范式进化投资组合'
$ws.Cells.Item(108, 1).Style = "Normal"
$ws.Cells.Item(108, 2).Formula = '''159781'
$ws.Cells.Item(108, 2).Style = "Normal"
$ws.Cells.Item(108, 3).Formula = '''科创创业ETF'
$ws.Cells.Item(108, 3).Style = "Normal"
$ws.Cells.Item(108, 4).Value = 5
$ws.Cells.Item(108, 5).Value = 9428.611488061997
$ws.Cells.Item(108, 6).Value = 0.53
$ws.Cells.Item(108, 7).Value = 4997.164088672858
$ws.Cells.Item(108, 8).Value = 99955.80999104423
$ws.Cells.Item(108, 9).Formula = '''202506171600'
$ws.Cells.Item(108, 9).Style = "Normal"
$ws.Cells.Item(109, 1).Formula = '''范式进化投资组合'
$ws.Cells.Item(109, 1).Style = "Normal"
$ws.Cells.Item(109, 2).Formula = '''510050'
$ws.Cells.Item(109, 2).Style = "Normal"
$ws.Cells.Item(109, 3).Formula = '''上证50ETF'
$ws.Cells.Item(109, 3).Style = "Normal"
$ws.Cells.Item(109, 4).Value = 5
$ws.Cells.Item(109, 5).Value = 1817.150577699221
$ws.Cells.Item(109, 6).Value = 2.75
$ws.Cells.Item(109, 7).Value = 4997.164088672857
$ws.Cells.Item(109, 8).Value = 99955.80999104423
$ws.Cells.Item(109, 9).Formula = '''202506171600'
$ws.Cells.Item(109, 9).Style = "Normal"
$ws.Cells.Item(110, 1).Formula = '''范式进化投资组合'
$ws.Cells.Item(110, 1).Style = "Normal"
$ws.Cells.Item(110, 2).Formula = '''510300'
$ws.Cells.Item(110, 2).Style = "Normal"
$ws.Cells.Item(110, 3).Formula = '''沪深300ETF'
$ws.Cells.Item(110, 3).Style = "Normal"
$ws.Cells.Item(110, 4).Value = 5
$ws.Cells.Item(110, 5).Value = 1252.422077361618
$ws.Cells.Item(110, 6).Value = 3.99
$ws.Cells.Item(110, 7).Value = 4997.164088672856
$ws.Cells.Item(110, 8).Value = 99955.80999104423
$ws.Cells.Item(110, 9).Formula = '''202506171600'
$ws.Cells.Item(110, 9).Style = "Normal"
$ws.Cells.Item(111, 1).Formula = '''范式进化投资组合'
$ws.Cells.Item(111, 1).Style = "Normal"
$ws.Cells.Item(111, 2).Formula = '''513100'
$ws.Cells.Item(111, 2).Style = "Normal"
$ws.Cells.Item(111, 3).Formula = '''纳指ETF'
$ws.Cells.Item(111, 3).Style = "Normal"
$ws.Cells.Item(111, 4).Value = 1
$ws.Cells.Item(111, 5).Value = 636.5814125697908
$ws.Cells.Item(111, 6).Value = 1.57
$ws.Cells.Item(111, 7).Value = 999.4328177345716
$ws.Cells.Item(111, 8).Value = 99955.80999104423
$ws.Cells.Item(111, 9).Formula = '''202506171600'
$ws.Cells.Item(111, 9).Style = "Normal"
$ws.Cells.Item(112, 1).Formula = '''范式进化投资组合'
$ws.Cells.Item(112, 1).Style = "Normal"
$ws.Cells.Item(112, 2).Formula = '''513290'
$ws.Cells.Item(112, 2).Style = "Normal"
$ws.Cells.Item(112, 3).Formula = '''纳指生物科技ETF'
$ws.Cells.Item(112, 3).Style = "Normal"
$ws.Cells.Item(112, 4).Value = 0.99
$ws.Cells.Item(112, 5).Value = 884.4538210040457
$ws.Cells.Item(112, 6).Value = 1.12
$ws.Cells.Item(112, 7).Value = 990.5882795245313
$ws.Cells.Item(112, 8).Value = 99955.80999104423
$ws.Cells.Item(112, 9).Formula = '''202506171600'
$ws.Cells.Item(112, 9).Style = "Normal"
$ws.Cells.Item(113, 1).Formula = '''范式进化投资组合'
$ws.Cells.Item(113, 1).Style = "Normal"
$ws.Cells.Item(113, 2).Formula = '''518880'
$ws.Cells.Item(113, 2).Style = "Normal"
$ws.Cells.Item(113, 3).Formula = '''黄金ETF'
$ws.Cells.Item(113, 3).Style = "Normal"
$ws.Cells.Item(113, 4).Value = 0.99
$ws.Cells.Item(113, 5).Value = 132.2001081659486
$ws.Cells.Item(113, 6).Value = 7.5
$ws.Cells.Item(113, 7).Value = 991.5008112446144
$ws.Cells.Item(113, 8).Value = 99955.80999104423
$ws.Cells.Item(113, 9).Formula = '''202506171600'
$ws.Cells.Item(113, 9).Style = "Normal"
$ws.Cells.Item(114, 1).Formula = '''范式进化投资组合'
$ws.Cells.Item(114, 1).Style = "Normal"
$ws.Cells.Item(114, 2).Formula = '''600085'
$ws.Cells.Item(114, 2).Style = "Normal"
$ws.Cells.Item(114, 3).Formula = '''同仁堂'
$ws.Cells.Item(114, 3).Style = "Normal"
$ws.Cells.Item(114, 4).Value = 1
$ws.Cells.Item(114, 5).Value = 27.49471300507762
$ws.Cells.Item(114, 6).Value = 36.19
$ws.Cells.Item(114, 7).Value = 995.033663653759
$ws.Cells.Item(114, 8).Value = 99955.80999104423
$ws.Cells.Item(114, 9).Formula = '''202506171600'
$ws.Cells.Item(114, 9).Style = "Normal"
$ws.Cells.Item(115, 1).Formula = '''范式进化投资组合'
$ws.Cells.Item(115, 1).Style = "Normal"
$ws.Cells.Item(115, 2).Formula = '''600900'
$ws.Cells.Item(115, 2).Style = "Normal"
$ws.Cells.Item(115, 3).Formula = '''长江电力'
$ws.Cells.Item(115, 3).Style = "Normal"
$ws.Cells.Item(115, 4).Value = 1.01
$ws.Cells.Item(115, 5).Value = 32.81132034584936
$ws.Cells.Item(115, 6).Value = 30.67
$ws.Cells.Item(115, 7).Value = 1006.3231950072
$ws.Cells.Item(115, 8).Value = 99955.80999104423
$ws.Cells.Item(115, 9).Formula = '''202506171600'
$ws.Cells.Item(115, 9).Style = "Normal"
$ws.Cells.Item(116, 1).Formula = '''范式进化投资组合'
$ws.Cells.Item(116, 1).Style = "Normal"
$ws.Cells.Item(116, 2).Formula = '''600989'
$ws.Cells.Item(116, 2).Style = "Normal"
$ws.Cells.Item(116, 3).Formula = '''宝丰能源'
$ws.Cells.Item(116, 3).Style = "Normal"
$ws.Cells.Item(116, 4).Value = 4.98
$ws.Cells.Item(116, 5).Value = 302.6749902285196
$ws.Cells.Item(116, 6).Value = 16.43
$ws.Cells.Item(116, 7).Value = 4972.950089454577
$ws.Cells.Item(116, 8).Value = 99955.80999104423
$ws.Cells.Item(116, 9).Formula = '''202506171600'
$ws.Cells.Item(116, 9).Style = "Normal"
$ws.Cells.Item(117, 1).Formula = '''范式进化投资组合'
$ws.Cells.Item(117, 1).Style = "Normal"
$ws.Cells.Item(117, 2).Formula = '''601899'
$ws.Cells.Item(117, 2).Style = "Normal"
$ws.Cells.Item(117, 3).Formula = '''XD紫金矿'
$ws.Cells.Item(117, 3).Style = "Normal"
$ws.Cells.Item(117, 4).Value = 10.03
$ws.Cells.Item(117, 5).Value = 531.0482559694855
$ws.Cells.Item(117, 6).Value = 18.87
$ws.Cells.Item(117, 7).Value = 10020.88059014419
$ws.Cells.Item(117, 8).Value = 99955.80999104423
$ws.Cells.Item(117, 9).Formula = '''202506171600'
$ws.Cells.Item(117, 9).Style = "Normal"
$ws.Cells.Item(118, 1).Formula = '''范式进化投资组合'
$ws.Cells.Item(118, 1).Style = "Normal"
$ws.Cells.Item(118, 2).Formula = '''603119'
$ws.Cells.Item(118, 2).Style = "Normal"
$ws.Cells.Item(118, 3).Formula = '''浙江荣泰'
$ws.Cells.Item(118, 3).Style = "Normal"
$ws.Cells.Item(118, 4).Value = 0.98
$ws.Cells.Item(118, 5).Value = 24.16423640557475
$ws.Cells.Item(118, 6).Value = 40.46
$ws.Cells.Item(118, 7).Value = 977.6850049695545
$ws.Cells.Item(118, 8).Value = 99955.80999104423
$ws.Cells.Item(118, 9).Formula = '''202506171600'
$ws.Cells.Item(118, 9).Style = "Normal"
$ws.Cells.Item(119, 1).Formula = '''范式进化投资组合'
$ws.Cells.Item(119, 1).Style = "Normal"
$ws.Cells.Item(119, 2).Formula = '''HK06881'
$ws.Cells.Item(119, 2).Style = "Normal"
$ws.Cells.Item(119, 3).Formula = '''中国银河'
$ws.Cells.Item(119, 3).Style = "Normal"
$ws.Cells.Item(119, 4).Value = 1
$ws.Cells.Item(119, 5).Value = 117.8576436007749
$ws.Cells.Item(119, 6).Value = 8.45
$ws.Cells.Item(119, 7).Value = 995.8970884265478
$ws.Cells.Item(119, 8).Value = 99955.80999104423
$ws.Cells.Item(119, 9).Formula = '''202506171600'
$ws.Cells.Item(119, 9).Style = "Normal"
$ws.Cells.Item(120, 1).Formula = '''范式进化投资组合'
$ws.Cells.Item(120, 1).Style = "Normal"
$ws.Cells.Item(120, 2).Formula = '''100000'
$ws.Cells.Item(120, 2).Style = "Normal"
$ws.Cells.Item(120, 3).Formula = '''现金'
$ws.Cells.Item(120, 3).Style = "Normal"
$ws.Cells.Item(120, 4).Value = 56.99
$ws.Cells.Item(120, 5).Value = 56967.67061087058
$ws.Cells.Item(120, 6).Value = 1
$ws.Cells.Item(120, 7).Value = 56967.67061087058
$ws.Cells.Item(120, 8).Value = 99955.80999104423
$ws.Cells.Item(120, 9).Formula = '''202506171600'
$ws.Cells.Item(120, 9).Style = "Normal"

